$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to text format so values like "1.00", "37.144.14" are preserved literally
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '37.144.14'
$ws.Range("E2").Value = '  +1.34%  '

# Row 3
$ws.Range("D3").Value = '2.055.16'
$ws.Range("E3").Value = '  -2.55%  '

# Row 4
$ws.Range("E4").Value = '  +0.03%  '

# Row 5
$ws.Range("D5").Value = '249.49'
$ws.Range("E5").Value = '  -1.57%  '

# Row 6
$ws.Range("D6").Value = '0.662'
$ws.Range("E6").Value = '  -0.84%  '

# Row 7
$ws.Range("E7").Value = '  +0.01%  '

# Row 8
$ws.Range("D8").Value = '56.37'
$ws.Range("E8").Value = '  +22.41%  '

# Row 9
$ws.Range("D9").Value = '62.27'
$ws.Range("E9").Value = '  -0.25%  '

# Row 10
$ws.Range("D10").Value = '0.383'
$ws.Range("E10").Value = '  +3.45%  '

# Row 11
$ws.Range("D11").Value = '0.0762'
$ws.Range("E11").Value = '  +2.85%  '

# Row 12
$ws.Range("E12").Value = '  +6.15%  '

# Row 13
$ws.Range("D13").Value = '15.20'
$ws.Range("E13").Value = '  +4.45%  '

# Row 14
$ws.Range("D14").Value = '2.359.42'
$ws.Range("E14").Value = '  -2.30%  '

# Row 15
$ws.Range("D15").Value = '0.833'
$ws.Range("E15").Value = '  -1.75%  '

# Row 16
$ws.Range("D16").Value = '5.29'
$ws.Range("E16").Value = '  +2.70%  '

# Row 17
$ws.Range("D17").Value = '2.061.31'
$ws.Range("E17").Value = '  -2.21%  '

# Row 18
$ws.Range("D18").Value = '37.066.95'
$ws.Range("E18").Value = '  +1.30%  '

# Row 19
$ws.Range("D19").Value = '72.79'
$ws.Range("E19").Value = '  -1.43%  '

# Row 20
$ws.Range("D20").Value = '14.59'
$ws.Range("E20").Value = '  +9.98%  '

# Row 21
$ws.Range("D21").Value = '0.0₃0861'
$ws.Range("E21").Value = '  +3.59%  '

# Row 22
$ws.Range("D22").Value = '238.61'
$ws.Range("E22").Value = '  -1.01%  '

# Row 23
$ws.Range("D23").Value = '5.28'
$ws.Range("E23").Value = '  +1.56%  '

# Row 24
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("D25").Value = '2.44'
$ws.Range("E25").Value = '  -0.57%  '

# Row 26
$ws.Range("D26").Value = '170.80'
$ws.Range("E26").Value = '  -0.33%  '

# Row 27
$ws.Range("D27").Value = '9.16'
$ws.Range("E27").Value = '  -0.29%  '

# Row 28
$ws.Range("D28").Value = '20.61'
$ws.Range("E28").Value = '  -4.77%  '

# Row 29
$ws.Range("D29").Value = '2.03'
$ws.Range("E29").Value = '  +1.67%  '

# Row 30
$ws.Range("D30").Value = '0.123'
$ws.Range("E30").Value = '  +0.09%  '

# Row 31
$ws.Range("E31").Value = '  +21.77%  '

# Row 32
$ws.Range("D32").Value = '22.70'
$ws.Range("E32").Value = '  -1.40%  '

# Row 33
$ws.Range("D33").Value = '4.58'
$ws.Range("E33").Value = '  +1.48%  '

# Row 34
$ws.Range("D34").Value = '0.0631'
$ws.Range("E34").Value = '  +4.69%  '

# Row 35
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").Value = '4.38'
$ws.Range("E35").Value = '  +6.24%  '

# Row 36
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").Value = '0.0901'
$ws.Range("E36").Value = '  -8.93%  '

# Row 37
$ws.Range("E37").Value = '  +0.17%  '

# Row 38
$ws.Range("D38").Value = '2.27'
$ws.Range("E38").Value = '  -3.36%  '

# Row 39
$ws.Range("D39").Value = '1.83'
$ws.Range("E39").Value = '  -3.15%  '

# Row 40
$ws.Range("D40").Value = '1.35'
$ws.Range("E40").Value = '  +0.71%  '

# Row 41
$ws.Range("E41").Value = '  +23.82%  '

# Row 42
$ws.Range("D42").Value = '18.05'
$ws.Range("E42").Value = '  +11.69%  '

# Row 43
$ws.Range("D43").Value = '0.0228'
$ws.Range("E43").Value = '  +3.27%  '

# Row 44
$ws.Range("D44").Value = '1.16'
$ws.Range("E44").Value = '  -3.08%  '

# Row 45
$ws.Range("D45").Value = '97.64'
$ws.Range("E45").Value = '  -1.60%  '

# Row 46
$ws.Range("B46").Value = 'HuobiToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D46").Value = '2.80'
$ws.Range("E46").Value = '  -1.90%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").Value = '2.43'
$ws.Range("E47").Value = '  +6.85%  '

# Row 48
$ws.Range("B48").Value = 'Maker'
$ws.Range("C48").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D48").Value = '1.311.03'
$ws.Range("E48").Value = '  -3.82%  '

# Row 49
$ws.Range("B49").Value = 'FTXToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D49").Value = '3.90'
$ws.Range("E49").Value = '  +40.29%  '

# Row 50
$ws.Range("D50").Value = '2.91'
$ws.Range("E50").Value = '  +3.14%  '

# Row 51
$ws.Range("D51").Value = '6.98'
$ws.Range("E51").Value = '  +4.18%  '
